$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert 3 new columns before column B ---------------------------------
# This shifts existing B,C,D,E (and their formatting) to E,F,G,H.
$ws.Columns("B:D").Insert()

# Keep the custom column width (8 characters) consistent across C:H, matching
# the width that previously applied to the single "UN" column.
$ws.Columns("C:H").ColumnWidth = 7.166666666666666

# --- 2. New date headers in row 1 (newest dates at the left) -----------------
$ws.Cells.Item(1, 2).Value = "Jun_27"
$ws.Cells.Item(1, 3).Value = "Jun_26"
$ws.Cells.Item(1, 4).Value = "Jun_26"

# --- 3. Fill the new B,C,D columns with the default "UN" rating for every ----
#        existing data row (2-27); row 22 gets the special downgrade note on
#        C22/D22 instead of the default value.
for ($r = 2; $r -le 27; $r++) {
    if ($r -eq 22) {
        $ws.Cells.Item($r, 2).Value = "UN"
        $ws.Cells.Item($r, 3).Value = "6/22/2018,Downgrades,Buy -> Hold,"
        $ws.Cells.Item($r, 4).Value = "6/22/2018,Downgrades,Buy -> Hold,"
    } else {
        $ws.Cells.Item($r, 2).Value = "UN"
        $ws.Cells.Item($r, 3).Value = "UN"
        $ws.Cells.Item($r, 4).Value = "UN"
    }
}

# --- 4. New analyst rows at the bottom of the table ---------------------------
$ws.Cells.Item(28, 1).Value = "Benchmark"
$ws.Cells.Item(28, 2).Value = "UN"
$ws.Cells.Item(28, 3).Value = "UN"
$ws.Cells.Item(28, 4).Value = "UN"

$ws.Cells.Item(29, 1).Value = "Evercore ISI"
$ws.Cells.Item(29, 2).Value = "UN"
$ws.Cells.Item(29, 3).Value = "UN"
$ws.Cells.Item(29, 4).Value = "UN"

Write-Host "Edit complete"
